$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.710.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.203.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.99"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +13.63%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.38%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.32"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.90%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.51%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.81%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.529.74"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.196.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.636.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.80"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.83%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.51"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.01%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.86%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.92"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0863"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.58%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0361"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.45%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.13%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.86%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.54"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.40%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.24%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.72%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.29"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.34%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0981"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.71%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.88%  "

# Row 50
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.432"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.40%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.63%  "
